$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (and restore the correct
# Coin/Link pairing for rows 45-48, which had drifted out of sync).
# A leading apostrophe forces Excel to keep the cell as text instead of
# reinterpreting a plain-looking decimal (e.g. "243.76") as a number.

$ws.Range("D2").Value = "29.152.14"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.838.56"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'243.76"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.6252"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D8").Value = "'0.07519"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "'0.2944"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'23.31"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "'0.07706"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "1.848.87"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "'5.022"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "'0.6772"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'83.22"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'0.000009290"
$ws.Range("E16").Value = "  -5.15%  "
$ws.Range("D17").Value = "'5.976"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "29.168.11"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "2.096.18"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'232.96"
$ws.Range("E20").Value = "  +2.70%  "
$ws.Range("D21").Value = "'12.71"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'7.180"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'160.62"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'0.1405"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'8.555"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'17.97"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'1.495"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'4.188"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'4.151"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").Value = "'0.05574"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'0.7528"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "'1.854"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "'1.147"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "'2.669"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "1.243.52"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'2.769"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.01791"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'6.610"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.9037"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'102.32"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.991.08"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.86"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5088"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("D49").Value = "'0.4090"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "'0.07417"
$ws.Range("E51").Value = "  +14.86%  "
